$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 0.6

$ws.Range("H12").Select()
